$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.730.90"
$ws.Range("E2").Value = "  +6.85%  "
$ws.Range("D3").Value = "1.812.66"
$ws.Range("E3").Value = "  +4.99%  "
$ws.Range("D4").Value = "'0.9992"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'250.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.57%  "
$ws.Range("D6").Value = "'0.9993"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "'0.4985"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.95%  "
$ws.Range("E8").Value = "  +7.32%  "
$ws.Range("D9").Value = "'0.06389"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.87%  "
$ws.Range("D10").Value = "1.808.32"
$ws.Range("E10").Value = "  +4.71%  "
$ws.Range("D11").Value = "'16.74"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.75%  "
$ws.Range("D12").Value = "'0.07130"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.32%  "
$ws.Range("D13").Value = "'0.6476"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.49%  "
$ws.Range("D14").Value = "'4.702"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.82%  "
$ws.Range("D15").Value = "'81.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.00%  "
$ws.Range("D16").Value = "28.697.35"
$ws.Range("E16").Value = "  +7.69%  "
$ws.Range("D17").Value = "'0.9993"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "'0.000007388"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.95%  "
$ws.Range("D19").Value = "'0.9989"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("E20").Value = "  +7.26%  "
$ws.Range("D21").Value = "2.040.17"
$ws.Range("E21").Value = "  +4.44%  "
$ws.Range("D22").Value = "'4.620"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.41%  "
$ws.Range("D23").Value = "'8.885"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.73%  "
$ws.Range("D24").Value = "'5.334"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.57%  "
$ws.Range("D25").Value = "'143.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.47%  "
$ws.Range("D26").Value = "'16.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.84%  "
$ws.Range("D27").Value = "'1.875"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.43%  "
$ws.Range("D28").Value = "'112.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.21%  "
$ws.Range("D29").Value = "'1.390"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("D30").Value = "'4.181"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.03%  "
$ws.Range("D31").Value = "'0.08351"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.46%  "
$ws.Range("D32").Value = "'3.840"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.13%  "
$ws.Range("D33").Value = "'0.04972"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +9.92%  "
$ws.Range("D34").Value = "'1.089"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.96%  "
$ws.Range("D35").Value = "'0.6758"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.15%  "
$ws.Range("D36").Value = "'2.668"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.76%  "
$ws.Range("D37").Value = "'2.727"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.02%  "
$ws.Range("D38").Value = "'0.9613"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.44%  "
$ws.Range("D39").Value = "'2.146"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.63%  "
$ws.Range("E40").Value = "  +5.78%  "
$ws.Range("D41").Value = "'5.979"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.18%  "
$ws.Range("D42").Value = "'0.9999"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'101.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.86%  "
$ws.Range("D44").Value = "'0.4111"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.70%  "
$ws.Range("D45").Value = "'7.194"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.64%  "
$ws.Range("E46").Value = "  +5.42%  "
$ws.Range("D47").Value = "'0.05498"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.76%  "
$ws.Range("D48").Value = "'8.160"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.75%  "
$ws.Range("D49").Value = "'31.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.33%  "
$ws.Range("D50").Value = "'0.3631"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.59%  "
$ws.Range("D51").Value = "'1.302"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.57%  "
